# Proper accounting of tax credits for CC retrofits after introduction of
# parallel routes for CO2 storage: add a "regime" column to the
# production_tax_credits sheet, tag the existing CC-retrofit / DAC rows as
# "OB3", and add a parallel set of rows with the "IRA 2022" credit values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("production_tax_credits")
$ws.Activate()

# --- New "regime" header column (J) -----------------------------------
$ws.Range("J1").Value = "regime"

# --- Tag the existing CC-retrofit / DAC rows (10-15) as "OB3" ---------
$ws.Range("J10").Value = "OB3"
$ws.Range("J11").Value = "OB3"
$ws.Range("J12").Value = "OB3"
$ws.Range("J13").Value = "OB3"
$ws.Range("J14").Value = "OB3"
$ws.Range("J15").Value = "OB3"

# --- Insert 6 new rows (16-21) for the "IRA 2022" regime values -------
# (this pushes the old rows 16-18 -- Alkaline electrolyzer large, PEM
# electrolyzer, SOEC -- down to rows 22-24)
$ws.Range("A16:A21").EntireRow.Insert()

$newRows = @(
    @{ Row = 16; Carrier = "ethanol from starch CC"; B = -60; C = 20; D = 12; E = 0.07 },
    @{ Row = 17; Carrier = "SMR CC";                 B = -60; C = 20; D = 12; E = 0.07 },
    @{ Row = 18; Carrier = "DRI CC";                 B = -60; C = 20; D = 12; E = 0.07 },
    @{ Row = 19; Carrier = "BF-BOF CC";               B = -60; C = 20; D = 12; E = 0.07 },
    @{ Row = 20; Carrier = "dry clinker CC";          B = -60; C = 20; D = 12; E = 0.07 },
    @{ Row = 21; Carrier = "DAC";                     B = -130; C = 30; D = 12; E = 0.1 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.Carrier
    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $r.D
    $ws.Range("E$row").Value = $r.E
    $ws.Range("F$row").Formula = "=E$row/(1-(1/(1+E$row)^C$row))"
    $ws.Range("G$row").Formula = "=F$row/(1-(1/(1+F$row)^D$row))"
    $ws.Range("I$row").Formula = "=B$row/(1-`$H`$2)*F$row/G$row"
    $ws.Range("J$row").Value = "IRA 2022"
}

# --- Selection / view bookkeeping (best effort) ------------------------
$ws.Range("I10:J21").Select()
$ws.Range("J21").Activate()

$ws2 = $wb.Worksheets.Item("investment_tax_credits")
$ws2.Range("I10:J21").Select()
$ws2.Range("H5").Activate()

$ws.Activate()
